# Generate Report for Handback
#
# Renames the two handed-back source files throughout the workbook:
#   ab6cf907-ef68-4f0f-9ae1-a47026b937b8  ->  cfd4e294-a26b-4277-828b-d295242f3449
#   d2f17490-2dad-4f74-9f05-6ee478d42efb  ->  ffffa777109d-9d96-4144-9506-70369ee9803c
# and refreshes the generated xlf hashes / timestamps that go with the new
# handback run. (Both files now hash to the same new content, so the
# zh-cn/de-de sheets' two rows converge on the same generated .xlf name.)

$wb = $excel.ActiveWorkbook

$oldId1 = "ab6cf907-ef68-4f0f-9ae1-a47026b937b8"
$newId1 = "cfd4e294-a26b-4277-828b-d295242f3449"
$oldId2 = "d2f17490-2dad-4f74-9f05-6ee478d42efb"
$newId2 = "ffffa777109d-9d96-4144-9506-70369ee9803c"

$oldHash1 = "83c057c332cad23ee05fca9d24080b06bc355d72"
$newHash = "9946fd911840eda06b637f3d5a4481b4a832470f"
$oldHash2 = "c94ed268f6c6903dd1d3c94885ba664e6db4a140"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newId1.md"
$ws1.Range("G2").Value = "2016-08-26 04:59:41"

$ws1.Range("A3").Value = "$newId2.md"
$ws1.Range("G3").Value = "2016-08-26 04:59:41"

# Rebuild the hyperlinks on this sheet with refreshed display text but the
# same target addresses (those did not change) and the same rId order.
$addr1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/$oldId1.md"
$addr2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/$oldId2.md"

$ws1.Hyperlinks.Delete()
$ws1.Range("B2").Value = "e2e\$newId1.md"
$ws1.Range("B3").Value = "e2e\$newId2.md"
[void]$ws1.Hyperlinks.Add($ws1.Range("B2"), $addr1, "", "", "e2e\$newId1.md")
[void]$ws1.Hyperlinks.Add($ws1.Range("B3"), $addr2, "", "", "e2e\$newId2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newId1.md"
$ws2.Range("G2").Value = "$newId1.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-26 04:59:37"
$ws2.Range("I2").Value = "$newId1.md"
$ws2.Range("J2").Value = "$newId1.$newHash.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-26 04:59:54"

$ws2.Range("A3").Value = "$newId2.md"
$ws2.Range("G3").Value = "$newId1.$newHash.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-26 04:59:37"
$ws2.Range("I3").Value = "$newId2.md"
$ws2.Range("J3").Value = "$newId1.$newHash.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-26 04:59:54"

$addr1zh = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/$oldId1.md"
$addr1zhcn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/299f355ceb56ad10487e14a8d42e3b27b1fd6fba/e2e/$oldId1.md"
$addr2zh = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/$oldId2.md"
$addr2zhcn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/299f355ceb56ad10487e14a8d42e3b27b1fd6fba/e2e/$oldId2.md"

$ws2.Hyperlinks.Delete()
[void]$ws2.Hyperlinks.Add($ws2.Range("A2"), $addr1zh, "", "", "$newId1.md")
[void]$ws2.Hyperlinks.Add($ws2.Range("I2"), $addr1zhcn, "", "", "$newId1.md")
[void]$ws2.Hyperlinks.Add($ws2.Range("A3"), $addr2zh, "", "", "$newId2.md")
[void]$ws2.Hyperlinks.Add($ws2.Range("I3"), $addr2zhcn, "", "", "$newId2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newId1.md"
$ws3.Range("G2").Value = "$newId1.$newHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-26 04:59:41"
$ws3.Range("I2").Value = "$newId1.md"
$ws3.Range("J2").Value = "$newId1.$newHash.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-26 05:00:03"

$ws3.Range("A3").Value = "$newId2.md"
$ws3.Range("G3").Value = "$newId1.$newHash.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-26 04:59:41"
$ws3.Range("I3").Value = "$newId2.md"
$ws3.Range("J3").Value = "$newId1.$newHash.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-26 05:00:03"

$addr1de = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/$oldId1.md"
$addr1dede = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5667c8244383cfb1d175130456b9e7ddc68a5bf7/e2e/$oldId1.md"
$addr2de = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce2e7d31840497bd5957871f63225c7151799e25/e2e/$oldId2.md"
$addr2dede = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5667c8244383cfb1d175130456b9e7ddc68a5bf7/e2e/$oldId2.md"

$ws3.Hyperlinks.Delete()
[void]$ws3.Hyperlinks.Add($ws3.Range("A2"), $addr1de, "", "", "$newId1.md")
[void]$ws3.Hyperlinks.Add($ws3.Range("I2"), $addr1dede, "", "", "$newId1.md")
[void]$ws3.Hyperlinks.Add($ws3.Range("A3"), $addr2de, "", "", "$newId2.md")
[void]$ws3.Hyperlinks.Add($ws3.Range("I3"), $addr2dede, "", "", "$newId2.md")
